$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# stay as text so they match the original inline-string cell type.
$textCells = @(
    'D5',
    'D6',
    'D9',
    'D13',
    'D19',
    'D20',
    'D24',
    'D29',
    'D30',
    'D36',
    'D37',
    'D38',
    'D44',
    'D46',
    'D48',
    'D49'
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '68.183.63'
$ws.Range('E2').Value = '  -1.13%  '
$ws.Range('D3').Value = '2.642.24'
$ws.Range('E3').Value = '  -1.02%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '597.61'
$ws.Range('E5').Value = '  -0.74%  '
$ws.Range('D6').Value = '156.14'
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -1.16%  '
$ws.Range('D9').Value = '0.141'
$ws.Range('E9').Value = '  +1.09%  '
$ws.Range('E10').Value = '  -1.41%  '
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('D13').Value = '27.99'
$ws.Range('E13').Value = '  -0.49%  '
$ws.Range('E14').Value = '  +0.29%  '
$ws.Range('D15').Value = '3.123.07'
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').Value = '68.218.28'
$ws.Range('E16').Value = '  -0.90%  '
$ws.Range('D17').Value = '2.635.27'
$ws.Range('E17').Value = '  -1.30%  '
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('D19').Value = '363.09'
$ws.Range('E19').Value = '  -1.34%  '
$ws.Range('D20').Value = '7.34'
$ws.Range('E20').Value = '  -1.73%  '
$ws.Range('E21').Value = '  +2.72%  '
$ws.Range('E22').Value = '  -2.65%  '
$ws.Range('E23').Value = '  -3.22%  '
$ws.Range('D24').Value = '75.18'
$ws.Range('E24').Value = '  +3.02%  '
$ws.Range('E25').Value = '  -0.39%  '
$ws.Range('E26').Value = '  -3.78%  '
$ws.Range('E27').Value = '  +7.02%  '
$ws.Range('D28').Value = '2.775.92'
$ws.Range('E28').Value = '  -1.17%  '
$ws.Range('D29').Value = '0.0000105'
$ws.Range('E29').Value = '  -2.14%  '
$ws.Range('D30').Value = '555.09'
$ws.Range('E30').Value = '  -4.82%  '
$ws.Range('E31').Value = '  +0.29%  '
$ws.Range('E32').Value = '  -1.36%  '
$ws.Range('E33').Value = '  -1.18%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('E35').Value = '  -2.70%  '
$ws.Range('D36').Value = '1.55'
$ws.Range('E36').Value = '  -0.30%  '
$ws.Range('D37').Value = '161.82'
$ws.Range('E37').Value = '  +1.89%  '
$ws.Range('D38').Value = '19.58'
$ws.Range('E38').Value = '  +1.23%  '
$ws.Range('E39').Value = '  +0.59%  '
$ws.Range('E40').Value = '  -3.70%  '
$ws.Range('E41').Value = '  -1.78%  '
$ws.Range('D42').Value = '0.0₆0335'
$ws.Range('E42').Value = '  +3.28%  '
$ws.Range('D44').Value = '2.61'
$ws.Range('E44').Value = '  -2.72%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').Value = '158.75'
$ws.Range('E46').Value = '  +1.17%  '
$ws.Range('E47').Value = '  -0.49%  '
$ws.Range('D48').Value = '22.08'
$ws.Range('D49').Value = '1.68'
$ws.Range('E49').Value = '  -2.77%  '
$ws.Range('E50').Value = '  -0.12%  '
$ws.Range('E51').Value = '  -1.11%  '

# Restore default (style-less) formatting on the forced-text cells so
# no extra style index remains attached, matching the source layout.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
